# Replace the "Games:" list (B1:B12) with a shorter "Locations" list (B1:B5),
# clearing the now-unused cells B6:B12 entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Locations"
$ws.Range("B2").Value = "L1"
$ws.Range("B3").Value = "L2"
$ws.Range("B4").Value = "L3"
$ws.Range("B5").Value = "L4"

# Remove the old "X vs Y" matchup cells that are no longer needed.
$ws.Range("B6:B12").Clear()

# Update the active selection to match the new editing position.
$ws.Range("B5").Select()
